$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header row (A1:C1) ---
# Old headers: "No." / "Code" / "NAME"  ->  New headers: "no" / "product_code" / "product_name"
$ws.Range("A1").Value = "no"
$ws.Range("B1").Value = "product_code"
$ws.Range("C1").Value = "product_name"

# --- Update column C data rows ---
# Old "combined" values ("V - SHIRT", "G - HANDBAG", "B - NECKLACE") are replaced by the
# plain product names already used elsewhere in the sheet ("V Shirt", "G Handbag", "B Necklace").
$ws.Range("C2").Value = "V Shirt"
$ws.Range("C3").Value = "G Handbag"
$ws.Range("C4").Value = "B Necklace"

# --- Remove the now unused D and E columns (split-code helper columns) ---
$ws.Range("D1:E4").Clear()

# --- Restore the selection to match the saved view state ---
[void]$ws.Range("C3").Select()
